# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" (fund-holdings detail for the new
#   quarter) positioned between the existing "总计" summary sheet and the
#   "2021-Q2" detail sheet.
# - Add the corresponding summary row ("2022-Q3": 4 holdings, 0.24 亿元)
#   to the "总计" sheet, ahead of the pre-existing "2021-Q2" summary row.

function Set-TextValue($range, $value) {
    # Force the written value to be stored as text (inline/shared string)
    # rather than being auto-coerced to a number, even when the text looks
    # numeric (e.g. "4.30", "0.1552", "630010") — then drop the temporary
    # text number-format so no stray style index is left on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2021-Q2")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet right before "2021-Q2" so the
#    tab order becomes: 总计, 2022-Q3, 2021-Q2.
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Header row (B1:H1) + index column (A2:A5) reuse the same bold/boxed
# style already used on the "总计" sheet's header/index cells.
$totalSheet.Range("B1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$q3Sheet.Range("A2:A5").PasteSpecial(-4122)

$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

$q3Sheet.Range("A2").Value = 0
Set-TextValue $q3Sheet.Range("B2") "630010"
Set-TextValue $q3Sheet.Range("C2") "华商价值精选混合"
Set-TextValue $q3Sheet.Range("D2") "4.30"
Set-TextValue $q3Sheet.Range("E2") "81.81"
Set-TextValue $q3Sheet.Range("F2") "3.61"
Set-TextValue $q3Sheet.Range("G2") "0.1552"
$q3Sheet.Range("H2").Value = 5

$q3Sheet.Range("A3").Value = 1
Set-TextValue $q3Sheet.Range("B3") "011686"
Set-TextValue $q3Sheet.Range("C3") "创金合信先进装备股票C"
Set-TextValue $q3Sheet.Range("D3") "0.57"
Set-TextValue $q3Sheet.Range("E3") "80.17"
Set-TextValue $q3Sheet.Range("F3") "6.50"
Set-TextValue $q3Sheet.Range("G3") "0.0370"
$q3Sheet.Range("H3").Value = 6

$q3Sheet.Range("A4").Value = 2
Set-TextValue $q3Sheet.Range("B4") "630006"
Set-TextValue $q3Sheet.Range("C4") "华商产业升级混合"
Set-TextValue $q3Sheet.Range("D4") "0.85"
Set-TextValue $q3Sheet.Range("E4") "81.97"
Set-TextValue $q3Sheet.Range("F4") "3.63"
Set-TextValue $q3Sheet.Range("G4") "0.0309"
$q3Sheet.Range("H4").Value = 5

$q3Sheet.Range("A5").Value = 3
Set-TextValue $q3Sheet.Range("B5") "011685"
Set-TextValue $q3Sheet.Range("C5") "创金合信先进装备股票A"
Set-TextValue $q3Sheet.Range("D5") "0.25"
Set-TextValue $q3Sheet.Range("E5") "80.17"
Set-TextValue $q3Sheet.Range("F5") "6.50"
Set-TextValue $q3Sheet.Range("G5") "0.0162"
$q3Sheet.Range("H5").Value = 6

# ---------------------------------------------------------------------
# 2. Insert a matching summary row on "总计" (pushes the existing
#    "2021-Q2" summary row from row 2 down to row 3), then fill both
#    rows' values.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Drop the formatting Excel auto-copied from row 1 into the new row 2 —
# in the target, B2:D2 carry no explicit style (like the original row).
$totalSheet.Range("B2:D2").ClearFormats()

# A2 should keep the same boxed/bold index style the data row always had;
# re-apply it (Insert left A2 unstyled) by copying down from A3 (the old
# A2, shifted down, which still carries it).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.24

$totalSheet.Range("A3").Value = 1
